$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing departure/return dates ---
# Re-typed with a leading apostrophe so Excel keeps them as text (quote-prefixed)
# instead of converting them into date serial numbers.
$ws.Range("C3").Value = "'2022-03-03"
$ws.Range("D3").Value = "'2022-03-27"

# --- New "Email" column (with a mailto hyperlink on the value cell) ---
$ws.Range("E2").Value = "Email"
$ws.Range("E3").Value = "milosmarkovic@test.com"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:milosmarkovic@test.com")

# --- New "FirstName" / "LastName" columns ---
$ws.Range("F2").Value = "FirstName"
$ws.Range("G2").Value = "LastName"

$ws.Range("F3").Value = "Milos"
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("G3").Value = "Markovic"

# --- New "Gender" column ---
$ws.Range("H2").Value = "Gender"
$ws.Range("H3").Value = "Female"

# --- Column widths for the new columns ---
$ws.Columns.Item(5).ColumnWidth = 28.498697916666668
$ws.Columns.Item(6).ColumnWidth = 13.166666666666666

# --- Row height for the data row ---
$ws.Rows.Item(3).RowHeight = 13

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection shown when the sheet is reopened ---
$ws.Range("J7").Select() | Out-Null

Write-Host "done"
